{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Summary of changes (see commit message / xml diff):\n//  1) In the first (\"AAPL\") Python-script code block, remove the\n//     \"from datetime import datetime\" line and the following\n//     \"import pandas as pd\" line (they are no longer imported there).\n//  2) In the second (\"AMZN\") Python-script code block:\n//       - \"from datetime import datetime\" -> \"from time import sleep\"\n//       - remove \"import pandas as pd\"\n//       - remove \"import os\" / \"import time\" / \"import sys\"\n//       - \"time.sleep(20)\" -> \"sleep(20)\"\n//  3) Insert two new paragraphs right before the\n//     \"dataset['Date'] = pandas.to_datetime(\" paragraph, explaining that\n//     Power BI already imports pandas/matplotlib and linking to a blog\n//     post about performance tips (with the _GoBack bookmark landing in\n//     the middle of the first new paragraph, matching the author's last\n//     edit position).\n//  4) Simplify the \"The DAX pages shows...\" paragraph into a single run\n//     (dropping the old _GoBack bookmark that used to sit inside it).\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Helper: find the first paragraph whose text matches `predicate`.\n// ---------------------------------------------------------------------\nasync function findParagraph(predicate) {\n  const paras = body.paragraphs;\n  paras.load(\"items/text\");\n  await context.sync();\n  for (const p of paras.items) {\n    if (predicate(p.text)) {\n      return p;\n    }\n  }\n  return null;\n}\n\n// ---------------------------------------------------------------------\n// 1) First code block: drop \"from datetime import datetime\" (+ its\n//    trailing whitespace run) and the following \"import pandas as pd\".\n// ---------------------------------------------------------------------\nconst firstDatetimeImport = await findParagraph(\n  (t) => t.indexOf(\"from datetime import datetime\") === 0 && t.trim() === \"from datetime import datetime\"\n);\nif (firstDatetimeImport) {\n  firstDatetimeImport.delete();\n  await context.sync();\n}\n\nconst firstPandasImport = await findParagraph((t) => t === \"import pandas as pd\");\nif (firstPandasImport) {\n  firstPandasImport.delete();\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Second code block: rewrite \"from datetime import datetime\" (exact\n//    match, no trailing spaces) into \"from time import sleep\".\n// ---------------------------------------------------------------------\nconst secondDatetimeImport = await findParagraph((t) => t === \"from datetime import datetime\");\nif (secondDatetimeImport) {\n  const r = secondDatetimeImport.getRange();\n  r.insertText(\"from time import sleep\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Remove \"import pandas as pd\" (second occurrence, now the only one left).\nconst secondPandasImport = await findParagraph((t) => t === \"import pandas as pd\");\nif (secondPandasImport) {\n  secondPandasImport.delete();\n  await context.sync();\n}\n\n// Remove \"import os\".\nconst importOs = await findParagraph((t) => t === \"import os\");\nif (importOs) {\n  importOs.delete();\n  await context.sync();\n}\n\n// Remove \"import time\".\nconst importTime = await findParagraph((t) => t === \"import time\");\nif (importTime) {\n  importTime.delete();\n  await context.sync();\n}\n\n// Remove \"import sys\".\nconst importSys = await findParagraph((t) => t === \"import sys\");\nif (importSys) {\n  importSys.delete();\n  await context.sync();\n}\n\n// \"time.sleep(20)\" -> \"sleep(20)\" (keep the italic+underline run formatting).\nconst timeSleepPara = await findParagraph((t) => t === \"time.sleep(20)\");\nif (timeSleepPara) {\n  const r = timeSleepPara.getRange();\n  r.insertText(\"sleep(20)\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Insert the two new paragraphs before\n//    \"dataset['Date'] = pandas.to_datetime(...\".\n// ---------------------------------------------------------------------\nconst pandasToDatetimePara = await findParagraph(\n  (t) => t.indexOf(\"dataset['Date'] = pandas.to_datetime(\") === 0\n);\nif (pandasToDatetimePara) {\n  const insertionPoint = pandasToDatetimePara.getRange(Word.RangeLocation.start);\n  insertionPoint.insertText(\n    \"If you open the query in you default Python IDE, you can see that Power BI already impo\" +\n      \"rts pandas as pandas and matplotlib.pyplot. Importing unnecessary libraries would slow down a script.\\r\" +\n      \"More reading about performance : https://dataveld.com/2018/11/10/5-performance-tips-for-r-and-python-scripts-in-power-bi/\\r\",\n    Word.InsertLocation.before\n  );\n  await context.sync();\n\n  // Place the \"_GoBack\" bookmark right after \"...already impo\", splitting\n  // that sentence into two runs exactly like the target document.\n  const splitResults = body.search(\"already impo\", { matchCase: true });\n  splitResults.load(\"items\");\n  await context.sync();\n  if (splitResults.items.length > 0) {\n    const splitPoint = splitResults.items[0].getRange(Word.RangeLocation.end);\n    splitPoint.insertBookmark(\"_GoBack\");\n    await context.sync();\n  }\n}\n\n// ---------------------------------------------------------------------\n// 4) Simplify the DAX paragraph into a single run; this also removes the\n//    old \"_GoBack\" bookmark that used to live inside it.\n// ---------------------------------------------------------------------\nconst daxPara = await findParagraph((t) => t.indexOf(\"The DAX pages shows\") === 0);\nif (daxPara) {\n  const r = daxPara.getRange();\n  r.insertText(\n    \"The DAX pages shows some basic DAX examples. You can see the DAX codes behind the measures by clicking on them.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Summary of changes (see commit message / xml diff):\n#  1) In the first (\"AAPL\") Python-script code block, remove the\n#     \"from datetime import datetime\" line and the following\n#     \"import pandas as pd\" line (they are no longer imported there).\n#  2) In the second (\"AMZN\") Python-script code block:\n#       - \"from datetime import datetime\" -> \"from time import sleep\"\n#       - remove \"import pandas as pd\"\n#       - remove \"import os\" / \"import time\" / \"import sys\"\n#       - \"time.sleep(20)\" -> \"sleep(20)\"\n#  3) Insert two new paragraphs right before the\n#     \"dataset['Date'] = pandas.to_datetime(\" paragraph, explaining that\n#     Power BI already imports pandas/matplotlib and linking to a blog\n#     post about performance tips (with the _GoBack bookmark landing in\n#     the middle of the first new paragraph, matching the author's last\n#     edit position).\n#  4) Simplify the \"The DAX pages shows...\" paragraph into a single run\n#     (dropping the old _GoBack bookmark that used to sit inside it).\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphByExactText($doc, [string]$text) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        $t = ($p.Range.Text -replace \"[\\r\\x0B]+$\", \"\")\n        if ($t -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Get-ParagraphByPrefix($doc, [string]$prefix) {\n    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {\n        $p = $doc.Paragraphs.Item($i)\n        if ($p.Range.Text.StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# ---------------------------------------------------------------------\n# 1) First code block: drop \"from datetime import datetime\" (+ its\n#    trailing whitespace run) and the following \"import pandas as pd\".\n# ---------------------------------------------------------------------\n$firstDatetimeImport = Get-ParagraphByExactText $d \"from datetime import datetime                               \"\nif ($firstDatetimeImport) {\n    $firstDatetimeImport.Range.Delete()\n}\n\n$firstPandasImport = Get-ParagraphByExactText $d \"import pandas as pd\"\nif ($firstPandasImport) {\n    $firstPandasImport.Range.Delete()\n}\n\n# ---------------------------------------------------------------------\n# 2) Second code block: rewrite \"from datetime import datetime\" (exact\n#    match, no trailing spaces) into \"from time import sleep\".\n# ---------------------------------------------------------------------\n$secondDatetimeImport = Get-ParagraphByExactText $d \"from datetime import datetime\"\nif ($secondDatetimeImport) {\n    $r = $secondDatetimeImport.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = \"from time import sleep\"\n}\n\n# Remove \"import pandas as pd\" (second occurrence, now the only one left).\n$secondPandasImport = Get-ParagraphByExactText $d \"import pandas as pd\"\nif ($secondPandasImport) {\n    $secondPandasImport.Range.Delete()\n}\n\n# Remove \"import os\".\n$importOs = Get-ParagraphByExactText $d \"import os\"\nif ($importOs) {\n    $importOs.Range.Delete()\n}\n\n# Remove \"import time\".\n$importTime = Get-ParagraphByExactText $d \"import time\"\nif ($importTime) {\n    $importTime.Range.Delete()\n}\n\n# Remove \"import sys\".\n$importSys = Get-ParagraphByExactText $d \"import sys\"\nif ($importSys) {\n    $importSys.Range.Delete()\n}\n\n# \"time.sleep(20)\" -> \"sleep(20)\" (keep the italic+underline run formatting).\n$timeSleepPara = Get-ParagraphByExactText $d \"time.sleep(20)\"\nif ($timeSleepPara) {\n    $r = $timeSleepPara.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Text = \"sleep(20)\"\n}\n\n# ---------------------------------------------------------------------\n# 3) Insert the two new paragraphs before\n#    \"dataset['Date'] = pandas.to_datetime(...\".\n# ---------------------------------------------------------------------\n$pandasToDatetimePara = Get-ParagraphByPrefix $d \"dataset['Date'] = pandas.to_datetime(\"\nif ($pandasToDatetimePara) {\n    $insertionPoint = $pandasToDatetimePara.Range\n    $insertionPoint.Collapse(1) | Out-Null  # wdCollapseStart\n    $insertionPoint.InsertBefore(\n        \"If you open the query in you default Python IDE, you can see that Power BI already impo\" +\n        \"rts pandas as pandas and matplotlib.pyplot. Importing unnecessary libraries would slow down a script.`r\" +\n        \"More reading about performance : https://dataveld.com/2018/11/10/5-performance-tips-for-r-and-python-scripts-in-power-bi/`r\"\n    )\n\n    # Place the \"_GoBack\" bookmark right after \"...already impo\", splitting\n    # that sentence into two runs exactly like the target document.\n    $splitRange = $d.Content\n    $found = $splitRange.Find.Execute(\"already impo\")\n    if ($found) {\n        $endPoint = $d.Range($splitRange.End, $splitRange.End)\n        $d.Bookmarks.Add(\"_GoBack\", $endPoint) | Out-Null\n    }\n}\n\n# ---------------------------------------------------------------------\n# 4) Simplify the DAX paragraph into a single run; this also removes the\n#    old \"_GoBack\" bookmark that used to live inside it.\n# ---------------------------------------------------------------------\n$daxPara = Get-ParagraphByPrefix $d \"The DAX pages shows\"\nif ($daxPara) {\n    $r = $daxPara.Range\n    $r.MoveEnd(1, -1) | Out-Null\n    $r.Delete()\n    $r.InsertAfter(\"The DAX pages shows some basic DAX examples. You can see the DAX codes behind the measures by clicking on them.\")\n}\n"}
